$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 97, shifting rows 97:142 down to 98:143
$ws.Rows("97:97").Insert()

# Populate the newly inserted row 97 with the new weekly entry
$ws.Range("A97").Value = 3
$ws.Range("B97").Value = "Femacal de La Calera"
$ws.Range("C97").Value = "Coquimbo"
$ws.Range("D97").Value = 44510
$ws.Range("E97").Value = 5
$ws.Range("F97").Value = 100112010
$ws.Range("G97").Value = "Achicoria"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 125
$ws.Range("K97").Value = 5500
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = 5760
$ws.Range("N97").Value = "$/caja 16 unidades"
$ws.Range("O97").Value = "Provincia de Quillota"
$ws.Range("P97").Value = 360
$ws.Range("Q97").Value = 16
$ws.Range("R97").Value = "Hortaliza"

# Match the date-number format style used by column D elsewhere
$ws.Range("D97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
